$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 8, duplicating row 7 (ANZ / Financials / 2024-03-01 data)
# but with Stock Code "Joshi     " instead of "ANZ       " -- this is the
# "Seconf rule check" row added based on the Threshold fix.
$ws.Range("A8").Value = "Joshi     "
$ws.Range("B8").Value = $ws.Range("B7").Value2
$ws.Range("C8").Value = $ws.Range("C7").Value2
$ws.Range("D8").Value = 28.9
$ws.Range("E8").Value = 29.15
$ws.Range("F8").Value = 950000
$ws.Range("G8").Value = $ws.Range("G7").Value2
$ws.Range("H8").Value = 82.1

# Selection moves to B2 after the edit
[void]$ws.Range("B2").Select()

# Column G (Trade Date) gets an explicit best-fit-style width once the new
# row is in place.
$ws.Columns("G:G").ColumnWidth = 10
